$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 40
$ws.Range("H40").Value2 = 1984
$ws.Range("I40").Value2 = 1478.2727
$ws.Range("J40").Value2 = 3374.75
$ws.Range("K40").Value2 = 1478.2727
$ws.Range("L40").Value2 = 3374.75
$ws.Range("M40").Value2 = -1303.2727
$ws.Range("N40").Value2 = -3724.75
# row 64
$ws.Range("H64").Value2 = 3521.9167
$ws.Range("I64").Value2 = 3272
$ws.Range("J64").Value2 = 3700.4285
$ws.Range("K64").Value2 = 3272
$ws.Range("L64").Value2 = 3700.4285
$ws.Range("M64").Value2 = -3024
$ws.Range("N64").Value2 = -4196.4285
# row 67
$ws.Range("H67").Value2 = 3521.9167
$ws.Range("I67").Value2 = 3272
$ws.Range("J67").Value2 = 3700.4285
$ws.Range("K67").Value2 = 3272
$ws.Range("L67").Value2 = 3700.4285
$ws.Range("M67").Value2 = -2414
$ws.Range("N67").Value2 = -5416.4285
# row 70
$ws.Range("H70").Value2 = 987.5
$ws.Range("I70").Value2 = 900
$ws.Range("K70").Value2 = 2700
$ws.Range("M70").Value2 = -2430
# row 73
$ws.Range("H73").Value2 = 987.5
$ws.Range("I73").Value2 = 900
$ws.Range("K73").Value2 = 2700
$ws.Range("M73").Value2 = -1764
# row 74
$ws.Range("H74").Value2 = 5817.3687
$ws.Range("I74").Value2 = 6948.1113
$ws.Range("J74").Value2 = 4799.7
$ws.Range("K74").Value2 = 6948.1113
$ws.Range("L74").Value2 = 4799.7
$ws.Range("M74").Value2 = -6012.1113
$ws.Range("N74").Value2 = -6671.7
# row 77
$ws.Range("H77").Value2 = 5817.3687
$ws.Range("I77").Value2 = 6948.1113
$ws.Range("J77").Value2 = 4799.7
$ws.Range("K77").Value2 = 34740.5565
$ws.Range("L77").Value2 = 23998.5
$ws.Range("M77").Value2 = -30060.5565
$ws.Range("N77").Value2 = -33358.5
# row 116
$ws.Range("H116").Value2 = 2557481.5
$ws.Range("I116").Value2 = 11906629
$ws.Range("J116").Value2 = 7713.909
$ws.Range("K116").Value2 = 11906629
$ws.Range("L116").Value2 = 7713.909
$ws.Range("M116").Value2 = -11903187
$ws.Range("N116").Value2 = -14597.909
# row 132
$ws.Range("H132").Value2 = 2526755.2
$ws.Range("I132").Value2 = 2778780.8
$ws.Range("J132").Value2 = 6500
$ws.Range("K132").Value2 = 8336342.399999999
$ws.Range("L132").Value2 = 19500
$ws.Range("M132").Value2 = -8333812.399999999
$ws.Range("N132").Value2 = -24560
# row 138
$ws.Range("H138").Value2 = 2302.705
$ws.Range("I138").Value2 = 797.8444
$ws.Range("J138").Value2 = 4354.788
$ws.Range("K138").Value2 = 2393.5332
$ws.Range("L138").Value2 = 13064.364
$ws.Range("M138").Value2 = 2746.4668
$ws.Range("N138").Value2 = -23344.364
# row 141
$ws.Range("H141").Value2 = 5193.3784
$ws.Range("I141").Value2 = 3683.3333
$ws.Range("J141").Value2 = 6222.9546
$ws.Range("K141").Value2 = 11049.9999
$ws.Range("L141").Value2 = 18668.8638
$ws.Range("M141").Value2 = -5869.999899999999
$ws.Range("N141").Value2 = -29028.8638

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 4
$ws.Range("H4").Value2 = 112.9
$ws.Range("I4").Value2 = 114.333336
$ws.Range("K4").Value2 = 114.333336
$ws.Range("M4").Value2 = 1.666663999999997
# row 9
$ws.Range("H9").Value2 = 79509
$ws.Range("J9").Value2 = 79509
$ws.Range("L9").Value2 = 79509
$ws.Range("N9").Value2 = -79849
# row 20
$ws.Range("H20").Value2 = 79509
$ws.Range("J20").Value2 = 79509
$ws.Range("L20").Value2 = 79509
$ws.Range("N20").Value2 = -80049
# row 80
$ws.Range("H80").Value2 = 27200
$ws.Range("J80").Value2 = 27200
$ws.Range("L80").Value2 = 27200
$ws.Range("N80").Value2 = -29196
# row 83
$ws.Range("H83").Value2 = 27200
$ws.Range("J83").Value2 = 27200
$ws.Range("L83").Value2 = 81600
$ws.Range("N83").Value2 = -91584
# row 97
$ws.Range("H97").Value2 = 559.9286
$ws.Range("I97").Value2 = 578
$ws.Range("J97").Value2 = 505.7143
$ws.Range("K97").Value2 = 578
$ws.Range("L97").Value2 = 505.7143
$ws.Range("M97").Value2 = -82
$ws.Range("N97").Value2 = -1497.7143

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value2 = 1543.2222
$ws.Range("I86").Value2 = 1523.625
$ws.Range("J86").Value2 = 1700
$ws.Range("K86").Value2 = 1523.625
$ws.Range("L86").Value2 = 1700
$ws.Range("M86").Value2 = -400.625
$ws.Range("N86").Value2 = -3946
# row 89
$ws.Range("H89").Value2 = 1543.2222
$ws.Range("I89").Value2 = 1523.625
$ws.Range("J89").Value2 = 1700
$ws.Range("K89").Value2 = 7618.125
$ws.Range("L89").Value2 = 8500
$ws.Range("M89").Value2 = -2002.125
$ws.Range("N89").Value2 = -19732
# row 94
$ws.Range("H94").Value2 = 592.1429000000001
$ws.Range("I94").Value2 = 565.8333
$ws.Range("K94").Value2 = 565.8333
$ws.Range("M94").Value2 = -114.8333
# row 134
$ws.Range("H134").Value2 = 2033.8889
$ws.Range("I134").Value2 = 1908.75
$ws.Range("J134").Value2 = 2471.875
$ws.Range("K134").Value2 = 5726.25
$ws.Range("L134").Value2 = 7415.625
$ws.Range("M134").Value2 = -3191.25
$ws.Range("N134").Value2 = -12485.625

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value2 = 336.22223
$ws.Range("J5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("N5").ClearContents()
# row 134
$ws.Range("H134").Value2 = 3614.5
$ws.Range("I134").Value2 = 1098.0476
$ws.Range("J134").Value2 = 4788.844
$ws.Range("K134").Value2 = 3294.142800000001
$ws.Range("L134").Value2 = 14366.532
$ws.Range("M134").Value2 = 1775.857199999999
$ws.Range("N134").Value2 = -24506.532
# row 135
$ws.Range("H135").Value2 = 336.22223
$ws.Range("J135").Value2 = 0
$ws.Range("L135").Value2 = 0
$ws.Range("N135").ClearContents()
# row 139
$ws.Range("H139").Value2 = 2597.138
$ws.Range("I139").Value2 = 1055
$ws.Range("J139").Value2 = 9999.4
$ws.Range("K139").Value2 = 3165
$ws.Range("L139").Value2 = 29998.2
$ws.Range("M139").Value2 = 1975
$ws.Range("N139").Value2 = -40278.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 43
$ws.Range("H43").Value2 = 12516.704
$ws.Range("I43").Value2 = 3825.1667
$ws.Range("J43").Value2 = 15000
$ws.Range("K43").Value2 = 3825.1667
$ws.Range("L43").Value2 = 15000
$ws.Range("M43").Value2 = -3674.1667
$ws.Range("N43").Value2 = -15302
# row 46
$ws.Range("H46").Value2 = 10010.25
$ws.Range("I46").Value2 = 10010.25
$ws.Range("J46").Value2 = 0
$ws.Range("K46").Value2 = 10010.25
$ws.Range("L46").Value2 = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value2 = -9854.25
# row 57
$ws.Range("H57").Value2 = 40000
$ws.Range("J57").Value2 = 40000
$ws.Range("L57").Value2 = 40000
$ws.Range("N57").Value2 = -41640
# row 80
$ws.Range("H80").Value2 = 1957.5714
$ws.Range("I80").Value2 = 1400
$ws.Range("J80").Value2 = 2267.3333
$ws.Range("K80").Value2 = 1400
$ws.Range("L80").Value2 = 2267.3333
$ws.Range("M80").Value2 = -402
$ws.Range("N80").Value2 = -4263.3333
# row 83
$ws.Range("H83").Value2 = 1957.5714
$ws.Range("I83").Value2 = 1400
$ws.Range("J83").Value2 = 2267.3333
$ws.Range("K83").Value2 = 7000
$ws.Range("L83").Value2 = 11336.6665
$ws.Range("M83").Value2 = -2008
$ws.Range("N83").Value2 = -21320.6665
# row 102
$ws.Range("H102").Value2 = 4993.2
$ws.Range("I102").Value2 = 4418.857
$ws.Range("J102").Value2 = 6333.3335
$ws.Range("K102").Value2 = 4418.857
$ws.Range("L102").Value2 = 6333.3335
$ws.Range("M102").Value2 = -2796.857
$ws.Range("N102").Value2 = -9577.333500000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 93
$ws.Range("H93").Value2 = 7784.294
$ws.Range("I93").Value2 = 34566.332
$ws.Range("K93").Value2 = 34566.332
$ws.Range("M93").Value2 = -33318.332
# row 132
$ws.Range("H132").Value2 = 3008.9644
$ws.Range("I132").Value2 = 1830.1538
$ws.Range("J132").Value2 = 4030.6
$ws.Range("K132").Value2 = 5490.4614
$ws.Range("L132").Value2 = 12091.8
$ws.Range("M132").Value2 = -2960.4614
$ws.Range("N132").Value2 = -17151.8
